$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (shifts existing rows 4-41 down to 5-42),
# then fill it with the new "Calories" nutrient entry.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = "Calories"
$ws.Cells.Item(4, 2).Value = 2400
$ws.Cells.Item(4, 3).Value = 2000

# Row 4 gets a slightly tighter height in the final layout.
$ws.Rows.Item(4).RowHeight = 13.8

# Move the active selection to C4 (matches the post-edit selection in the file).
$ws.Range("C4").Select() | Out-Null
